$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

# Row 2
Set-TextValue "D2" '29.103.47'
Set-TextValue "E2" '  +0.14%  '

# Row 3
Set-TextValue "D3" '1.838.23'
Set-TextValue "E3" '  +0.10%  '

# Row 4
Set-TextValue "D4" '0.9974'
Set-TextValue "E4" '  -0.31%  '

# Row 5
Set-TextValue "D5" '243.32'
Set-TextValue "E5" '  -0.30%  '

# Row 6
Set-TextValue "D6" '0.6245'
Set-TextValue "E6" '  -1.27%  '

# Row 7
Set-TextValue "D7" '0.9992'
Set-TextValue "E7" '  -0.17%  '

# Row 8
Set-TextValue "D8" '0.07510'
Set-TextValue "E8" '  -0.92%  '

# Row 9
Set-TextValue "D9" '0.2947'
Set-TextValue "E9" '  -0.05%  '

# Row 10
Set-TextValue "D10" '23.37'
Set-TextValue "E10" '  +2.32%  '

# Row 11
Set-TextValue "D11" '0.07695'
Set-TextValue "E11" '  -0.66%  '

# Row 12
Set-TextValue "D12" '1.833.66'
Set-TextValue "E12" '  -0.51%  '

# Row 13
Set-TextValue "D13" '5.019'
Set-TextValue "E13" '  +0.31%  '

# Row 14
Set-TextValue "D14" '0.6772'
Set-TextValue "E14" '  +0.88%  '

# Row 15
Set-TextValue "D15" '83.09'
Set-TextValue "E15" '  -0.21%  '

# Row 16
Set-TextValue "D16" '0.000009370'
Set-TextValue "E16" '  -4.12%  '

# Row 17
Set-TextValue "D17" '5.982'
Set-TextValue "E17" '  -2.30%  '

# Row 18
Set-TextValue "D18" '29.102.83'
Set-TextValue "E18" '  +0.05%  '

# Row 19
Set-TextValue "D19" '2.087.34'
Set-TextValue "E19" '  +0.20%  '

# Row 20
Set-TextValue "D20" '12.69'
Set-TextValue "E20" '  +1.02%  '

# Row 21
Set-TextValue "D21" '225.74'
Set-TextValue "E21" '  -0.65%  '

# Row 22
Set-TextValue "D22" '0.9997'
Set-TextValue "E22" '  -0.04%  '

# Row 23
Set-TextValue "D23" '7.170'
Set-TextValue "E23" '  -1.13%  '

# Row 24
Set-TextValue "D24" '0.9989'
Set-TextValue "E24" '  -0.21%  '

# Row 25
Set-TextValue "D25" '160.14'

# Row 26
Set-TextValue "E26" '  -0.49%  '

# Row 27
Set-TextValue "D27" '8.544'
Set-TextValue "E27" '  -0.06%  '

# Row 28
Set-TextValue "D28" '17.92'
Set-TextValue "E28" '  -0.34%  '

# Row 29
Set-TextValue "D29" '1.494'
Set-TextValue "E29" '  -0.63%  '

# Row 30
Set-TextValue "D30" '4.189'
Set-TextValue "E30" '  +1.50%  '

# Row 31
Set-TextValue "B31" 'InternetComputer(DFINITY)'
Set-TextValue "C31" 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue "D31" '4.150'
Set-TextValue "E31" '  +2.25%  '

# Row 32
Set-TextValue "B32" 'Hedera'
Set-TextValue "C32" 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue "D32" '0.05590'
Set-TextValue "E32" '  +4.21%  '

# Row 34
Set-TextValue "D34" '0.7500'
Set-TextValue "E34" '  -0.15%  '

# Row 35
Set-TextValue "E35" '  -0.71%  '

# Row 36
Set-TextValue "D36" '1.145'
Set-TextValue "E36" '  +0.22%  '

# Row 37
Set-TextValue "D37" '2.662'
Set-TextValue "E37" '  -0.25%  '

# Row 38
Set-TextValue "D38" '1.238.59'
Set-TextValue "E38" '  -1.01%  '

# Row 39
Set-TextValue "D39" '2.770'
Set-TextValue "E39" '  +0.32%  '

# Row 40
Set-TextValue "D40" '0.01783'
Set-TextValue "E40" '  -0.86%  '

# Row 41
Set-TextValue "D41" '6.565'
Set-TextValue "E41" '  -0.37%  '

# Row 42
Set-TextValue "E42" '  -0.65%  '

# Row 43
Set-TextValue "D43" '0.9992'

# Row 44
Set-TextValue "D44" '102.57'
Set-TextValue "E44" '  -0.15%  '

# Row 45
Set-TextValue "D45" '1.984.87'
Set-TextValue "E45" '  +0.18%  '

# Row 46
Set-TextValue "D46" '66.55'
Set-TextValue "E46" '  +2.35%  '

# Row 47
Set-TextValue "D47" '0.00000000124'
Set-TextValue "E47" '  +0.17%  '

# Row 48
Set-TextValue "D48" '0.5080'
Set-TextValue "E48" '  -0.69%  '

# Row 49
Set-TextValue "D49" '0.4082'
Set-TextValue "E49" '  -0.28%  '

# Row 50
Set-TextValue "D50" '9.075'
Set-TextValue "E50" '  +0.21%  '

# Row 51
Set-TextValue "D51" '0.05842'
Set-TextValue "E51" '  +0.63%  '
